$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 16.44711831065678
$ws.Cells.Item(2, 3).Value = 8.651268738831185
$ws.Cells.Item(2, 4).Value = 7.254635730637227
$ws.Cells.Item(2, 6).Value = 47.84801966962823
$ws.Cells.Item(2, 7).Value = 3.733221758670223
$ws.Cells.Item(2, 10).Value = 11.57496939932736
$ws.Cells.Item(2, 11).Value = 15.96434949061308
$ws.Cells.Item(2, 14).Value = 23.04671829507067

$ws.Cells.Item(3, 2).Value = 16.24902804071727
$ws.Cells.Item(3, 3).Value = 8.509496279951062
$ws.Cells.Item(3, 4).Value = 7.22124987853009
$ws.Cells.Item(3, 6).Value = 47.7110484607768
$ws.Cells.Item(3, 7).Value = 3.736493423664342
$ws.Cells.Item(3, 10).Value = 11.56359162661096
$ws.Cells.Item(3, 11).Value = 15.8366100363417
$ws.Cells.Item(3, 14).Value = 23.08566712772199

$ws.Cells.Item(4, 2).Value = 16.1310178350153
$ws.Cells.Item(4, 3).Value = 8.424069023183931
$ws.Cells.Item(4, 4).Value = 7.202605359040039
$ws.Cells.Item(4, 6).Value = 47.63712113760913
$ws.Cells.Item(4, 7).Value = 3.738606189661548
$ws.Cells.Item(4, 10).Value = 11.55901292489441
$ws.Cells.Item(4, 11).Value = 15.76201489755569
$ws.Cells.Item(4, 14).Value = 23.11146454301359

$ws.Cells.Item(5, 2).Value = 16.08389489598204
$ws.Cells.Item(5, 3).Value = 8.389713224718376
$ws.Cells.Item(5, 4).Value = 7.195480533003608
$ws.Cells.Item(5, 6).Value = 47.609567204489
$ws.Cells.Item(5, 7).Value = 3.739493394451422
$ws.Cells.Item(5, 10).Value = 11.55775332308992
$ws.Cells.Item(5, 11).Value = 15.73261093989191
$ws.Cells.Item(5, 14).Value = 23.12245036278911

$ws.Cells.Item(6, 2).Value = 16.07613021861112
$ws.Cells.Item(6, 3).Value = 8.384037452081456
$ws.Cells.Item(6, 4).Value = 7.194326207660498
$ws.Cells.Item(6, 6).Value = 47.60514759139367
$ws.Cells.Item(6, 7).Value = 3.739642301518803
$ws.Cells.Item(6, 10).Value = 11.55758080264753
$ws.Cells.Item(6, 11).Value = 15.72778929596073
$ws.Cells.Item(6, 14).Value = 23.12430312055606

$ws.Cells.Item(7, 2).Value = 16.13037832829432
$ws.Cells.Item(7, 3).Value = 8.423603776658711
$ws.Cells.Item(7, 4).Value = 7.202507347859602
$ws.Cells.Item(7, 6).Value = 47.63673910542646
$ws.Cells.Item(7, 7).Value = 3.738618048463002
$ws.Cells.Item(7, 10).Value = 11.55899348183831
$ws.Cells.Item(7, 11).Value = 15.76161428304067
$ws.Cells.Item(7, 14).Value = 23.11161078619437

$ws.Cells.Item(8, 2).Value = 16.3781009170129
$ws.Cells.Item(8, 3).Value = 8.602078455551791
$ws.Cells.Item(8, 4).Value = 7.242743561556784
$ws.Cells.Item(8, 6).Value = 47.79868956441871
$ws.Cells.Item(8, 7).Value = 3.734328311928561
$ws.Cells.Item(8, 10).Value = 11.57054737712771
$ws.Cells.Item(8, 11).Value = 15.91952658679152
$ws.Cells.Item(8, 14).Value = 23.05975695160178

$ws.Cells.Item(9, 2).Value = 16.88997913640606
$ws.Cells.Item(9, 3).Value = 8.96277089154122
$ws.Cells.Item(9, 4).Value = 7.336039315361141
$ws.Cells.Item(9, 6).Value = 48.19627430273592
$ws.Cells.Item(9, 7).Value = 3.726736556981216
$ws.Cells.Item(9, 10).Value = 11.61224913495867
$ws.Cells.Item(9, 11).Value = 16.25826753044123
$ws.Cells.Item(9, 14).Value = 22.97302143090049

$ws.Cells.Item(10, 2).Value = 17.27837636401071
$ws.Cells.Item(10, 3).Value = 9.23140849071569
$ws.Cells.Item(10, 4).Value = 7.412884260596453
$ws.Cells.Item(10, 6).Value = 48.53601987356416
$ws.Cells.Item(10, 7).Value = 3.721652882042298
$ws.Cells.Item(10, 10).Value = 11.65439923789185
$ws.Cells.Item(10, 11).Value = 16.52291626750192
$ws.Cells.Item(10, 14).Value = 22.91842685972998

$ws.Cells.Item(11, 2).Value = 17.45696965740889
$ws.Cells.Item(11, 3).Value = 9.353809834920272
$ws.Cells.Item(11, 4).Value = 7.449532734792828
$ws.Cells.Item(11, 6).Value = 48.70062933047283
$ws.Cells.Item(11, 7).Value = 3.719446137705424
$ws.Cells.Item(11, 10).Value = 11.67604473322306
$ws.Cells.Item(11, 11).Value = 16.6462961047479
$ws.Cells.Item(11, 14).Value = 22.89557667143087

$ws.Cells.Item(12, 2).Value = 17.52480837598991
$ws.Cells.Item(12, 3).Value = 9.40014015776301
$ws.Cells.Item(12, 4).Value = 7.463643535214404
$ws.Cells.Item(12, 6).Value = 48.76437995780131
$ws.Cells.Item(12, 7).Value = 3.718625620789275
$ws.Cells.Item(12, 10).Value = 11.68459346961234
$ws.Cells.Item(12, 11).Value = 16.69340839699266
$ws.Cells.Item(12, 14).Value = 22.88720973429297

$ws.Cells.Item(13, 2).Value = 17.5101898147941
$ws.Cells.Item(13, 3).Value = 9.390163758198387
$ws.Cells.Item(13, 4).Value = 7.460594337114777
$ws.Cells.Item(13, 6).Value = 48.75058761282227
$ws.Cells.Item(13, 7).Value = 3.718801662355035
$ws.Cells.Item(13, 10).Value = 11.68273674528983
$ws.Cells.Item(13, 11).Value = 16.6832451560746
$ws.Cells.Item(13, 14).Value = 22.88899897930941

$ws.Cells.Item(14, 2).Value = 17.46254696544358
$ws.Cells.Item(14, 3).Value = 9.357622140223569
$ws.Cells.Item(14, 4).Value = 7.450689033820831
$ws.Cells.Item(14, 6).Value = 48.70584592035276
$ws.Cells.Item(14, 7).Value = 3.719378330591425
$ws.Cells.Item(14, 10).Value = 11.67674100750205
$ws.Cells.Item(14, 11).Value = 16.65016444942155
$ws.Cells.Item(14, 14).Value = 22.89488258521424

$ws.Cells.Item(15, 2).Value = 17.43338968073731
$ws.Cells.Item(15, 3).Value = 9.33768538600596
$ws.Cells.Item(15, 4).Value = 7.444651753768472
$ws.Cells.Item(15, 6).Value = 48.67862390477084
$ws.Cells.Item(15, 7).Value = 3.719733524118927
$ws.Cells.Item(15, 10).Value = 11.67311418835972
$ws.Cells.Item(15, 11).Value = 16.62995132039156
$ws.Cells.Item(15, 14).Value = 22.8985237148286

$ws.Cells.Item(16, 2).Value = 17.26673769889068
$ws.Cells.Item(16, 3).Value = 9.223409031976088
$ws.Cells.Item(16, 4).Value = 7.410522351922208
$ws.Cells.Item(16, 6).Value = 48.52546215945457
$ws.Cells.Item(16, 7).Value = 3.721799220046525
$ws.Cells.Item(16, 10).Value = 11.65303410810859
$ws.Cells.Item(16, 11).Value = 16.51490995526213
$ws.Cells.Item(16, 14).Value = 22.91996015959183

$ws.Cells.Item(17, 2).Value = 17.16494449876148
$ws.Cells.Item(17, 3).Value = 9.15331939993329
$ws.Cells.Item(17, 4).Value = 7.3900110034155
$ws.Cells.Item(17, 6).Value = 48.43405800556793
$ws.Cells.Item(17, 7).Value = 3.723093502047758
$ws.Cells.Item(17, 10).Value = 11.64134634632315
$ws.Cells.Item(17, 11).Value = 16.44507408025606
$ws.Cells.Item(17, 14).Value = 22.93361950068524

$ws.Cells.Item(18, 2).Value = 17.1065793399078
$ws.Cells.Item(18, 3).Value = 9.113027032166606
$ws.Cells.Item(18, 4).Value = 7.378373264918983
$ws.Cells.Item(18, 6).Value = 48.38243334150432
$ws.Cells.Item(18, 7).Value = 3.723847906864745
$ws.Cells.Item(18, 10).Value = 11.63485667980404
$ws.Cells.Item(18, 11).Value = 16.40519012069856
$ws.Cells.Item(18, 14).Value = 22.94166281583715

$ws.Cells.Item(19, 2).Value = 17.08685143367579
$ws.Cells.Item(19, 3).Value = 9.099389904672258
$ws.Cells.Item(19, 4).Value = 7.374460689926565
$ws.Cells.Item(19, 6).Value = 48.36511790295517
$ws.Cells.Item(19, 7).Value = 3.724105050175656
$ws.Cells.Item(19, 10).Value = 11.63269947192851
$ws.Cells.Item(19, 11).Value = 16.39173601367181
$ws.Cells.Item(19, 14).Value = 22.94441821708421

$ws.Cells.Item(20, 2).Value = 17.1757620323817
$ws.Cells.Item(20, 3).Value = 9.160778711125205
$ws.Cells.Item(20, 4).Value = 7.392178001881383
$ws.Cells.Item(20, 6).Value = 48.44369016629184
$ws.Cells.Item(20, 7).Value = 3.722954692529603
$ws.Cells.Item(20, 10).Value = 11.64256645464365
$ws.Cells.Item(20, 11).Value = 16.45247913737948
$ws.Cells.Item(20, 14).Value = 22.93214610175796

$ws.Cells.Item(21, 2).Value = 17.4765356603454
$ws.Cells.Item(21, 3).Value = 9.367181335583677
$ws.Cells.Item(21, 4).Value = 7.453592227777382
$ws.Cells.Item(21, 6).Value = 48.71894944177908
$ws.Cells.Item(21, 7).Value = 3.719208539271574
$ws.Cells.Item(21, 10).Value = 11.67849257460642
$ws.Cells.Item(21, 11).Value = 16.65987075527073
$ws.Cells.Item(21, 14).Value = 22.89314666357847

$ws.Cells.Item(22, 2).Value = 17.67429907032666
$ws.Cells.Item(22, 3).Value = 9.501937121459253
$ws.Cells.Item(22, 4).Value = 7.495081915124607
$ws.Cells.Item(22, 6).Value = 48.90708894217534
$ws.Cells.Item(22, 7).Value = 3.716848354572828
$ws.Cells.Item(22, 10).Value = 11.70402261392076
$ws.Cells.Item(22, 11).Value = 16.79767398649097
$ws.Cells.Item(22, 14).Value = 22.86932514128454

$ws.Cells.Item(23, 2).Value = 17.56866141172149
$ws.Cells.Item(23, 3).Value = 9.430043880297916
$ws.Cells.Item(23, 4).Value = 7.472817924391667
$ws.Cells.Item(23, 6).Value = 48.80593162994919
$ws.Cells.Item(23, 7).Value = 3.718099994456359
$ws.Cells.Item(23, 10).Value = 11.69021034839574
$ws.Cells.Item(23, 11).Value = 16.72393192444851
$ws.Cells.Item(23, 14).Value = 22.88188645909106

$ws.Cells.Item(24, 2).Value = 17.17087093107272
$ws.Cells.Item(24, 3).Value = 9.15740634317115
$ws.Cells.Item(24, 4).Value = 7.391197819825342
$ws.Cells.Item(24, 6).Value = 48.43933258347678
$ws.Cells.Item(24, 7).Value = 3.72301741620643
$ws.Cells.Item(24, 10).Value = 11.64201412769294
$ws.Cells.Item(24, 11).Value = 16.44913048073389
$ws.Cells.Item(24, 14).Value = 22.93281163239564

$ws.Cells.Item(25, 2).Value = 16.74907626461129
$ws.Cells.Item(25, 3).Value = 8.864348412029678
$ws.Cells.Item(25, 4).Value = 7.309306071381232
$ws.Cells.Item(25, 6).Value = 48.08025874950734
$ws.Cells.Item(25, 7).Value = 3.728703135606488
$ws.Cells.Item(25, 10).Value = 11.59893734558593
$ws.Cells.Item(25, 11).Value = 16.16371219076167
$ws.Cells.Item(25, 14).Value = 22.99488358583801
